# Hindalco prices update: add a new "21-12-2025" row at the top of the
# daily price table (row 2), pushing every existing row down by one.
# The new row re-uses the Basic Price / Circular Date / Circular Link
# that were already the most recent ones (no fresh circular was issued
# for 21-12-2025 yet), exactly matching the prior latest row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new blank row at row 2 -------------------------------
# This shifts sheetData rows (values + cell styles) down by one row, and
# Excel keeps the worksheet dimension in sync automatically.
$ws.Rows.Item(2).Insert()

# --- Step 2: populate the new row 2 ----------------------------------------
# B..F should equal what used to be row 2 (now shifted to row 3); only the
# date in column A is new.
$ws.Range("A2").Value = "21-12-2025"
$ws.Range("B2").Value = $ws.Range("B3").Value()
$ws.Range("C2").Value = $ws.Range("C3").Value()
$ws.Range("D2").Value = $ws.Range("D3").Value()
$ws.Range("E2").Value = $ws.Range("E3").Value()
$ws.Range("F2").Value = $ws.Range("F3").Value()

# Make sure row 2 carries the same formatting (number format / alignment)
# as the rest of the data rows - copy formats only from row 3.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: rebuild hyperlinks ---------------------------------------------
# Row insertion does not relocate the worksheet's hyperlink anchors, so
# drop every hyperlink and re-create them from the (now-shifted) column F
# text, which always equals the link target in this sheet.
$ws.Hyperlinks.Delete()

$lastRow = $ws.UsedRange().Rows().Count()
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $link = $cell.Value()
    if ($link -ne $null -and $link.ToString().StartsWith("http")) {
        $ws.Hyperlinks.Add($cell, $link, "", "", $link)
    }
}
